$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.1633928571428571
$ws.Range("A3").Value = 0.5102678571428572
$ws.Range("A4").Value = 0.7857142857142858
$ws.Range("A5").Value = 0.9178571428571429
$ws.Range("A6").Value = 0.9651785714285714
$ws.Range("A7").Value = 0.9897321428571428
$ws.Range("A8").Value = 0.9977678571428571
$ws.Range("A9").Value = 0.9986607142857142
$ws.Range("A25").Value = 0.9991071428571429
$ws.Range("A26").Value = 0.9991071428571429
$ws.Range("A27").Value = 0.9991071428571429
